$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07033851334821055
$ws.Range("B3").Value = -0.4340269950258345
$ws.Range("B4").Value = 13.271069403387685

$ws.Rows.Item(5).Delete()
